$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the F1 header from "Path Length" to "Run Time"
$ws.Cells.Item(1, 6).Value = "Run Time"

# Data for rows 2-46: columns are (RowNumber, A=PuzzleNumber, D=SolutionLength, E=SearchPathLength, F=RunTime)
# Column B (Algorithm) and C (Heuristic) are unchanged by this edit.
$data = @(
    @(2, 1, 7, 286, 0.2240033149719238),
    @(3, 1, 7, 184, 0.1329855918884277),
    @(4, 1, 7, 251, 0.2160248756408691),
    @(5, 1, 7, 130, 0.1239974498748779),
    @(6, 1, 7, 171, 0.117988109588623),
    @(7, 1, 7, 20, 0.01198911666870117),
    @(8, 1, 9, 37, 0.02399539947509766),
    @(9, 1, 7, 20, 0.009990215301513672),
    @(10, 1, 8, 24, 0.01405882835388184),
    @(11, 2, 7, 736, 1.422996282577515),
    @(12, 2, 7, 397, 0.6589710712432861),
    @(13, 2, 7, 428, 0.7860367298126221),
    @(14, 2, 8, 203, 0.2419888973236084),
    @(15, 2, 8, 389, 0.6329941749572754),
    @(16, 2, 10, 27, 0.0260012149810791),
    @(17, 2, 8, 25, 0.01803398132324219),
    @(18, 2, 10, 27, 0.0230410099029541),
    @(19, 2, 9, 21, 0.0149993896484375),
    @(20, 3, 9, 246, 0.1779863834381104),
    @(21, 3, 9, 157, 0.0919651985168457),
    @(22, 3, 9, 178, 0.1279988288879395),
    @(23, 3, 10, 151, 0.09299373626708984),
    @(24, 3, 9, 152, 0.08903145790100098),
    @(25, 3, 15, 71, 0.03003764152526855),
    @(26, 3, 12, 82, 0.03300333023071289),
    @(27, 3, 15, 71, 0.03700828552246094),
    @(28, 3, 14, 74, 0.04003620147705078),
    @(29, 4, 9, 2552, 18.49196171760559),
    @(30, 4, 9, 1206, 5.473022222518921),
    @(31, 4, 9, 1513, 7.963034629821777),
    @(32, 4, 9, 669, 2.217032670974731),
    @(33, 4, 9, 1165, 5.524951219558716),
    @(34, 4, 9, 35, 0.0429990291595459),
    @(35, 4, 13, 84, 0.2020695209503174),
    @(36, 4, 9, 35, 0.04423403739929199),
    @(37, 4, 12, 34, 0.1520001888275146),
    @(38, 5, 5, 59, 0.03000378608703613),
    @(39, 5, 5, 51, 0.03101778030395508),
    @(40, 5, 5, 56, 0.02999639511108398),
    @(41, 5, 5, 48, 0.02102851867675781),
    @(42, 5, 5, 51, 0.02498412132263184),
    @(43, 5, 6, 20, 0.01099753379821777),
    @(44, 5, 6, 23, 0.01000857353210449),
    @(45, 5, 6, 20, 0.01403331756591797),
    @(46, 5, 6, 20, 0.01599526405334473)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]   # A: Puzzle Number
    $ws.Cells.Item($r, 4).Value = $row[2]   # D: Solution Length
    $ws.Cells.Item($r, 5).Value = $row[3]   # E: Search Path Length
    $ws.Cells.Item($r, 6).Value = $row[4]   # F: Run Time
}

# Remove the now-obsolete rows 47-55 (old extra puzzle entries)
$ws.Range("A47:F55").EntireRow.Delete()
